$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data added one more weekly record to the top of this block
# (row 195), pushing the existing rows 195-205 down to 196-206 while
# keeping their values unchanged. Insert a new row at 195 to reproduce
# that shift, then populate it with the new record.
$ws.Rows(195).Insert()

$ws.Cells.Item(195, 1).Value = 11
$ws.Cells.Item(195, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(195, 3).Value = "Bíobío"
$ws.Cells.Item(195, 4).Value = 44610
$ws.Cells.Item(195, 5).Value = 8
$ws.Cells.Item(195, 6).Value = 100112017
$ws.Cells.Item(195, 7).Value = "Apio"
$ws.Cells.Item(195, 8).Value = "Americana (o)"
$ws.Cells.Item(195, 9).Value = "Primera"
$ws.Cells.Item(195, 10).Value = 100
$ws.Cells.Item(195, 11).Value = 7000
$ws.Cells.Item(195, 12).Value = 8000
$ws.Cells.Item(195, 13).Value = 7500
$ws.Cells.Item(195, 14).Value = "$/docena de matas"
$ws.Cells.Item(195, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(195, 16).Value = 1250
$ws.Cells.Item(195, 17).Value = 6
$ws.Cells.Item(195, 18).Value = "Hortaliza"
